# The workbook tracks weekly price records for "Membrillo" (quince) at the
# Feria Lagunitas de Puerto Montt. A new weekly record (dated 2022-07-29,
# i.e. Excel serial 44771) was added to the data set. In the source system
# the records are ordered, so the new record lands right after the existing
# row 34 (2021-07-09) and everything that used to be row 35 onward shifts
# down by one row - growing the table from 96 to 97 data-bearing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; Excel shifts rows 35:96 down to 36:97
# and copies formatting (e.g. the date-number-format on column D) from the
# row above, matching how the row was actually authored.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value = "Los Lagos"
$ws.Cells.Item(35, 4).Value = 44771
$ws.Cells.Item(35, 5).Value = 10
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100104
$ws.Cells.Item(35, 8).Value = "Frutos de pepita"
$ws.Cells.Item(35, 9).Value = 100104003
$ws.Cells.Item(35, 10).Value = "Membrillo"
$ws.Cells.Item(35, 11).Value = "Champion"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 500
$ws.Cells.Item(35, 14).Value = 13500
$ws.Cells.Item(35, 15).Value = 14000
$ws.Cells.Item(35, 16).Value = 13750
$ws.Cells.Item(35, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(35, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(35, 19).Value = 764
$ws.Cells.Item(35, 20).Value = 18
